$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")
$ws.Activate()

# Write cells in the same order the strings were originally added to the
# shared string table so new shared-string indices line up with the source
# workbook (Clear 50 / Clear 600 / Clear 250, then the rest in row order).
$ws.Cells.Item(59, 2).Value = "Clear 50 chansons différentes"
$ws.Cells.Item(61, 2).Value = "Clear 600 chanson différentes"
$ws.Cells.Item(60, 2).Value = "Clear 250 chansons différentes"
$ws.Cells.Item(62, 2).Value = "Battez le boss 1"
$ws.Cells.Item(63, 2).Value = "Battez le boss 2"
$ws.Cells.Item(64, 2).Value = "Battez le boss 3"
$ws.Cells.Item(65, 2).Value = "Battez le boss 4"
$ws.Cells.Item(66, 2).Value = "Battez le boss 5"
$ws.Cells.Item(67, 2).Value = "Battez le boss 6"
$ws.Cells.Item(68, 2).Value = "Battez le boss 7"
$ws.Cells.Item(69, 2).Value = "Battez le boss 8"
$ws.Cells.Item(70, 2).Value = "Finissez Cublast"
$ws.Cells.Item(71, 2).Value = "Obtenez toutes les médailles de bronze dans le mode Story"
$ws.Cells.Item(72, 2).Value = "Obtenez toutes les médailles d'argent dans le mode Story"
$ws.Cells.Item(73, 2).Value = "Obtenez toutes les médailles d'or dans le mode Story"
$ws.Cells.Item(74, 2).Value = "Obtenez toutes les médailles de Quad dans le mode Story"

$ws.Range("B75").Select() | Out-Null

# Scroll the window so row 58 is the top visible row (mirrors the author's
# viewport move down to the newly appended rows).
$win = $excel.ActiveWindow
$win.ScrollRow = 58
$win.ScrollColumn = 1

